$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

Set-TextCell 'D2' '28.046.62'
$ws.Range('E2').Value = '  -0.61%  '
Set-TextCell 'D3' '1.901.35'
$ws.Range('E3').Value = '  +1.60%  '
Set-TextCell 'D4' '1.001'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextCell 'D5' '312.64'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('E6').Value = '  +0.05%  '
Set-TextCell 'D7' '0.5080'
$ws.Range('E7').Value = '  +0.58%  '
Set-TextCell 'D8' '0.3925'
$ws.Range('E8').Value = '  +0.22%  '
Set-TextCell 'D9' '0.09272'
$ws.Range('E9').Value = '  -3.68%  '
$ws.Range('E10').Value = '  -0.42%  '
Set-TextCell 'D11' '41.79'
$ws.Range('E11').Value = '  +2.32%  '
Set-TextCell 'D12' '6.366'
$ws.Range('E12').Value = '  -1.96%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D13' '1.902.14'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 'D14' '20.77'
$ws.Range('E14').Value = '  -0.82%  '
$ws.Range('E15').Value = '  +0.05%  '
Set-TextCell 'D16' '7.291'
$ws.Range('E16').Value = '  -1.94%  '
Set-TextCell 'D17' '0.00001117'
$ws.Range('E17').Value = '  -0.75%  '
Set-TextCell 'D18' '92.43'
$ws.Range('E18').Value = '  -0.54%  '
Set-TextCell 'D19' '0.06582'
$ws.Range('E19').Value = '  -0.68%  '
Set-TextCell 'D20' '17.78'
$ws.Range('E20').Value = '  +1.31%  '
Set-TextCell 'D22' '6.220'
$ws.Range('E22').Value = '  +1.14%  '
Set-TextCell 'D23' '28.097.80'
$ws.Range('E23').Value = '  -0.64%  '
Set-TextCell 'D24' '11.36'
$ws.Range('E24').Value = '  +0.18%  '
Set-TextCell 'D25' '2.322'
$ws.Range('E25').Value = '  +1.78%  '
Set-TextCell 'D26' '2.598'
$ws.Range('E26').Value = '  +2.36%  '
Set-TextCell 'D27' '2.117.22'
$ws.Range('E27').Value = '  +0.68%  '
Set-TextCell 'D28' '20.91'
$ws.Range('E28').Value = '  -1.36%  '
Set-TextCell 'D29' '157.48'
$ws.Range('E29').Value = '  +0.01%  '
Set-TextCell 'D30' '127.22'
$ws.Range('E30').Value = '  -0.04%  '
Set-TextCell 'D31' '1.084'
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('E32').Value = '  +0.92%  '
Set-TextCell 'D33' '5.605'
$ws.Range('E33').Value = '  -0.44%  '
Set-TextCell 'D34' '3.612'
$ws.Range('E34').Value = '  -0.36%  '
Set-TextCell 'D35' '9.584'
$ws.Range('E35').Value = '  -0.64%  '
Set-TextCell 'D36' '0.06654'
$ws.Range('E36').Value = '  -1.47%  '
Set-TextCell 'D37' '0.02407'
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D38' '1.226'
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextCell 'D39' '0.2172'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').Value = '  +6.72%  '
Set-TextCell 'D41' '0.6343'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D42' '11.42'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D43' '4.972'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('E44').Value = '  +0.01%  '
Set-TextCell 'D45' '13.27'
$ws.Range('E45').Value = '  -1.92%  '
$ws.Range('E46').Value = '  -0.77%  '
Set-TextCell 'D47' '3.706'
$ws.Range('E47').Value = '  +1.18%  '
$ws.Range('E48').Value = '  +0.69%  '
Set-TextCell 'D49' '2.008'
$ws.Range('E49').Value = '  +0.91%  '
Set-TextCell 'D50' '122.39'
$ws.Range('E50').Value = '  -1.87%  '
Set-TextCell 'D51' '1.178'
$ws.Range('E51').Value = '  -1.63%  '
